$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell I1 from "forecast_rule_evaluation_memory"
# to "forecast_rule_evaluation_memory_length"
$ws.Range("I1").Value = "forecast_rule_evaluation_memory_length"

# Update the active selection on the sheet from I10 to I6
$ws.Range("I6").Select()
